$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Re-arrange the "frames" / "docWidgets" / "docs" mini-tables on Sheet1.
#
# Column D ("frames" table): everything shifts up one row (the table used to
# start at D3, now starts at D2) and a new "userType" row is inserted right
# before "docWidgetID" (so the frames struct now also carries a userType
# field, matching docWidgets).
#
# Column F ("docWidgets" table): the old "userType" entry (F8) is dropped
# from here (frames now owns it) and "docID" moves up to take its place.
#
# Column H ("docs" table): shifts up one row to follow the column F change.
# ---------------------------------------------------------------------------

# Snapshot the old values first so later writes don't clobber reads.
$d3 = $ws.Range("D3").Value()
$d4 = $ws.Range("D4").Value()
$d5 = $ws.Range("D5").Value()

$f9 = $ws.Range("F9").Value()

$h8 = $ws.Range("H8").Value()
$h9 = $ws.Range("H9").Value()
$h10 = $ws.Range("H10").Value()
$h11 = $ws.Range("H11").Value()

# Column D: frames(3->2), ptr(4->3), attach(5->4), new userType at 5.
$ws.Range("D2").Value = $d3
$ws.Range("D3").Value = $d4
$ws.Range("D4").Value = $d5
$ws.Range("D5").Value = "userType"

# Column F: userType(8) removed, docID(9) moves up to 8.
$ws.Range("F8").Value = $f9
$ws.Range("F9").Value = ""

# Column H: docs/ID/ptr/name shift up one row (8..11 -> 7..10).
$ws.Range("H7").Value = $h8
$ws.Range("H8").Value = $h9
$ws.Range("H9").Value = $h10
$ws.Range("H10").Value = $h11
$ws.Range("H11").Value = ""

# Column D widened to fit the longer "mainWindowID" best-fit width.
$ws.Columns(4).ColumnWidth = 14.44140625

# Nudge the connector arrow that used to point at the old F8/H8 area so it
# tracks the row-8 (now row-7) content it was anchored to.
$shp = $ws.Shapes.Item(3)
$shp.Left = 302.8363779527559
$shp.Top = 106.55456692913386
$shp.Width = 97.1863779527559
$shp.Height = 0.8046456692913386

# Selection left where the author's cursor ended up.
$ws.Range("C10").Select() | Out-Null
